$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2208.7856  # H70: 6581090.5 -> 2208.7856
$ws.Cells.Item(70, 9).Value = 806  # I70: 20834122 -> 806
$ws.Cells.Item(70, 10).Value = 2988.111  # J70: 2768.4614 -> 2988.111
$ws.Cells.Item(70, 11).Value = 2418  # K70: 62502366 -> 2418
$ws.Cells.Item(70, 12).Value = 8964.332999999999  # L70: 8305.3842 -> 8964.332999999999
$ws.Cells.Item(70, 13).Value = -2148  # M70: -62502096 -> -2148
$ws.Cells.Item(70, 14).Value = -9504.332999999999  # N70: -8845.3842 -> -9504.332999999999
$ws.Cells.Item(73, 8).Value = 2208.7856  # H73: 6581090.5 -> 2208.7856
$ws.Cells.Item(73, 9).Value = 806  # I73: 20834122 -> 806
$ws.Cells.Item(73, 10).Value = 2988.111  # J73: 2768.4614 -> 2988.111
$ws.Cells.Item(73, 11).Value = 2418  # K73: 62502366 -> 2418
$ws.Cells.Item(73, 12).Value = 8964.332999999999  # L73: 8305.3842 -> 8964.332999999999
$ws.Cells.Item(73, 13).Value = -1482  # M73: -62501430 -> -1482
$ws.Cells.Item(73, 14).Value = -10836.333  # N73: -10177.3842 -> -10836.333
$ws.Cells.Item(80, 8).Value = 688.7037  # H80: 1014.75 -> 688.7037
$ws.Cells.Item(80, 9).Value = 541.2308  # I80: 617.3333 -> 541.2308
$ws.Cells.Item(80, 10).Value = 825.6429000000001  # J80: 1253.2 -> 825.6429000000001
$ws.Cells.Item(80, 11).Value = 1623.6924  # K80: 1851.9999 -> 1623.6924
$ws.Cells.Item(80, 12).Value = 2476.9287  # L80: 3759.6 -> 2476.9287
$ws.Cells.Item(80, 13).Value = -625.6924000000001  # M80: -853.9999 -> -625.6924000000001
$ws.Cells.Item(80, 14).Value = -4472.9287  # N80: -5755.6 -> -4472.9287
$ws.Cells.Item(83, 8).Value = 688.7037  # H83: 1014.75 -> 688.7037
$ws.Cells.Item(83, 9).Value = 541.2308  # I83: 617.3333 -> 541.2308
$ws.Cells.Item(83, 10).Value = 825.6429000000001  # J83: 1253.2 -> 825.6429000000001
$ws.Cells.Item(83, 11).Value = 4871.077200000001  # K83: 5555.9997 -> 4871.077200000001
$ws.Cells.Item(83, 12).Value = 7430.7861  # L83: 11278.8 -> 7430.7861
$ws.Cells.Item(83, 13).Value = 120.9227999999994  # M83: -563.9997000000003 -> 120.9227999999994
$ws.Cells.Item(83, 14).Value = -17414.7861  # N83: -21262.8 -> -17414.7861
$ws.Cells.Item(92, 8).Value = 881.8182  # H92: 1075 -> 881.8182
$ws.Cells.Item(92, 9).Value = 536.875  # I92: 599.1667 -> 536.875
$ws.Cells.Item(92, 10).Value = 1801.6666  # J92: 2502.5 -> 1801.6666
$ws.Cells.Item(92, 11).Value = 536.875  # K92: 599.1667 -> 536.875
$ws.Cells.Item(92, 12).Value = 1801.6666  # L92: 2502.5 -> 1801.6666
$ws.Cells.Item(92, 13).Value = 711.125  # M92: 648.8333 -> 711.125
$ws.Cells.Item(92, 14).Value = -4297.6666  # N92: -4998.5 -> -4297.6666
$ws.Cells.Item(113, 8).Value = 2879.3057  # H113: 2756.1667 -> 2879.3057
$ws.Cells.Item(113, 9).Value = 2836.7896  # I113: 2533.2222 -> 2836.7896
$ws.Cells.Item(113, 10).Value = 2926.8235  # J113: 3157.4666 -> 2926.8235
$ws.Cells.Item(113, 11).Value = 2836.7896  # K113: 2533.2222 -> 2836.7896
$ws.Cells.Item(113, 12).Value = 2926.8235  # L113: 3157.4666 -> 2926.8235
$ws.Cells.Item(113, 13).Value = 417.2103999999999  # M113: 720.7777999999998 -> 417.2103999999999
$ws.Cells.Item(113, 14).Value = -9434.8235  # N113: -9665.4666 -> -9434.8235
$ws.Cells.Item(132, 8).Value = 877267.4  # H132: 893245.6 -> 877267.4
$ws.Cells.Item(132, 9).Value = 2267.6667  # I132: 2432.8718 -> 2267.6667
$ws.Cells.Item(132, 10).Value = 3502266.5  # J132: 3064601.5 -> 3502266.5
$ws.Cells.Item(132, 11).Value = 6803.000100000001  # K132: 7298.6154 -> 6803.000100000001
$ws.Cells.Item(132, 12).Value = 10506799.5  # L132: 9193804.5 -> 10506799.5
$ws.Cells.Item(132, 13).Value = -4273.000100000001  # M132: -4768.6154 -> -4273.000100000001
$ws.Cells.Item(132, 14).Value = -10511859.5  # N132: -9198864.5 -> -10511859.5
$ws.Cells.Item(137, 8).Value = 2001502.9  # H137: 2129245 -> 2001502.9
$ws.Cells.Item(137, 9).Value = 3031367.2  # I137: 3226932 -> 3031367.2
$ws.Cells.Item(137, 10).Value = 2354.4707  # J137: 2476.5625 -> 2354.4707
$ws.Cells.Item(137, 11).Value = 9094101.600000001  # K137: 9680796 -> 9094101.600000001
$ws.Cells.Item(137, 12).Value = 7063.4121  # L137: 7429.6875 -> 7063.4121
$ws.Cells.Item(137, 13).Value = -9091551.600000001  # M137: -9678246 -> -9091551.600000001
$ws.Cells.Item(137, 14).Value = -12163.4121  # N137: -12529.6875 -> -12163.4121
$ws.Cells.Item(141, 8).Value = 1743.6482  # H141: 1432.7 -> 1743.6482
$ws.Cells.Item(141, 9).Value = 1177.4902  # I141: 1066.5 -> 1177.4902
$ws.Cells.Item(141, 10).Value = 11368.333  # J141: 12052.5 -> 11368.333
$ws.Cells.Item(141, 11).Value = 3532.4706  # K141: 3199.5 -> 3532.4706
$ws.Cells.Item(141, 12).Value = 34104.999  # L141: 36157.5 -> 34104.999
$ws.Cells.Item(141, 13).Value = 1647.5294  # M141: 1980.5 -> 1647.5294
$ws.Cells.Item(141, 14).Value = -44464.999  # N141: -46517.5 -> -44464.999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 26369750  # H61: 20876152 -> 26369750
$ws.Cells.Item(61, 9).Value = 29442352  # I61: 22751016 -> 29442352
$ws.Cells.Item(61, 11).Value = 29442352  # K61: 22751016 -> 29442352
$ws.Cells.Item(61, 13).Value = -29442140  # M61: -22750804 -> -29442140
$ws.Cells.Item(109, 8).Value = 35193  # H109: 52928.332 -> 35193
$ws.Cells.Item(109, 10).Value = 35193  # J109: 52928.332 -> 35193
$ws.Cells.Item(109, 12).Value = 35193  # L109: 52928.332 -> 35193
$ws.Cells.Item(109, 14).Value = -37967  # N109: -55702.332 -> -37967
$ws.Cells.Item(110, 8).Value = 556414.8  # H110: 435516.47 -> 556414.8
$ws.Cells.Item(110, 9).Value = 625766.7  # I110: 476841.84 -> 625766.7
$ws.Cells.Item(110, 11).Value = 625766.7  # K110: 476841.84 -> 625766.7
$ws.Cells.Item(110, 13).Value = -623721.7  # M110: -474796.84 -> -623721.7
$ws.Cells.Item(126, 8).Value = 4444  # H126: 0 -> 4444
$ws.Cells.Item(126, 9).Value = 4444  # I126: 0 -> 4444
$ws.Cells.Item(126, 11).Value = 13332  # K126: 0 -> 13332
$ws.Cells.Item(126, 13).Value = -10862  # M126: None -> -10862
$ws.Cells.Item(132, 8).Value = 54782.668  # H132: 50813.43 -> 54782.668
$ws.Cells.Item(132, 9).Value = 44378.477  # I132: 36440.895 -> 44378.477
$ws.Cells.Item(132, 10).Value = 69738.69  # J132: 79558.5 -> 69738.69
$ws.Cells.Item(132, 11).Value = 133135.431  # K132: 109322.685 -> 133135.431
$ws.Cells.Item(132, 12).Value = 209216.07  # L132: 238675.5 -> 209216.07
$ws.Cells.Item(132, 13).Value = -130605.431  # M132: -106792.685 -> -130605.431
$ws.Cells.Item(132, 14).Value = -214276.07  # N132: -243735.5 -> -214276.07
$ws.Cells.Item(136, 8).Value = 26369750  # H136: 20876152 -> 26369750
$ws.Cells.Item(136, 9).Value = 29442352  # I136: 22751016 -> 29442352
$ws.Cells.Item(136, 11).Value = 88327056  # K136: 68253048 -> 88327056
$ws.Cells.Item(136, 13).Value = -88324506  # M136: -68250498 -> -88324506

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(96, 8).Value = 13500  # H96: 25000 -> 13500
$ws.Cells.Item(96, 9).Value = 2000  # I96: 25000 -> 2000
$ws.Cells.Item(96, 11).Value = 2000  # K96: 25000 -> 2000
$ws.Cells.Item(96, 13).Value = 746  # M96: -22254 -> 746
$ws.Cells.Item(99, 8).Value = 984.63635  # H99: 930.6429000000001 -> 984.63635
$ws.Cells.Item(99, 9).Value = 925  # I99: 825.8 -> 925
$ws.Cells.Item(99, 10).Value = 1056.2  # J99: 1192.75 -> 1056.2
$ws.Cells.Item(99, 11).Value = 925  # K99: 825.8 -> 925
$ws.Cells.Item(99, 12).Value = 1056.2  # L99: 1192.75 -> 1056.2
$ws.Cells.Item(99, 13).Value = 573  # M99: 672.2 -> 573
$ws.Cells.Item(99, 14).Value = -4052.2  # N99: -4188.75 -> -4052.2
$ws.Cells.Item(107, 8).Value = 2342.7856  # H107: 2138.8125 -> 2342.7856
$ws.Cells.Item(107, 9).Value = 2066.5833  # I107: 1872.9286 -> 2066.5833
$ws.Cells.Item(107, 11).Value = 2066.5833  # K107: 1872.9286 -> 2066.5833
$ws.Cells.Item(107, 13).Value = -146.5832999999998  # M107: 47.07140000000004 -> -146.5832999999998
$ws.Cells.Item(128, 8).Value = 920  # H128: 1152.75 -> 920
$ws.Cells.Item(128, 9).Value = 920  # I128: 1152.75 -> 920
$ws.Cells.Item(128, 11).Value = 2760  # K128: 3458.25 -> 2760
$ws.Cells.Item(128, 13).Value = -270  # M128: -968.25 -> -270

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 3348000  # H4: 0 -> 3348000
$ws.Cells.Item(4, 10).Value = 3348000  # J4: 0 -> 3348000
$ws.Cells.Item(4, 12).Value = 3348000  # L4: 0 -> 3348000
$ws.Cells.Item(4, 14).Value = -3348224  # N4: None -> -3348224
$ws.Cells.Item(31, 8).Value = 3441  # H31: 3931.1 -> 3441
$ws.Cells.Item(31, 9).Value = 1386.6666  # I31: 1691.1875 -> 1386.6666
$ws.Cells.Item(31, 10).Value = 6759.5386  # J31: 6491 -> 6759.5386
$ws.Cells.Item(31, 11).Value = 1386.6666  # K31: 1691.1875 -> 1386.6666
$ws.Cells.Item(31, 12).Value = 6759.5386  # L31: 6491 -> 6759.5386
$ws.Cells.Item(31, 13).Value = -1091.6666  # M31: -1396.1875 -> -1091.6666
$ws.Cells.Item(31, 14).Value = -7349.5386  # N31: -7081 -> -7349.5386
$ws.Cells.Item(34, 8).Value = 3441  # H34: 3931.1 -> 3441
$ws.Cells.Item(34, 9).Value = 1386.6666  # I34: 1691.1875 -> 1386.6666
$ws.Cells.Item(34, 10).Value = 6759.5386  # J34: 6491 -> 6759.5386
$ws.Cells.Item(34, 11).Value = 1386.6666  # K34: 1691.1875 -> 1386.6666
$ws.Cells.Item(34, 12).Value = 6759.5386  # L34: 6491 -> 6759.5386
$ws.Cells.Item(34, 13).Value = -1184.6666  # M34: -1489.1875 -> -1184.6666
$ws.Cells.Item(34, 14).Value = -7163.5386  # N34: -6895 -> -7163.5386
$ws.Cells.Item(58, 8).Value = 27779928  # H58: 27029192 -> 27779928
$ws.Cells.Item(58, 9).Value = 37039052  # I58: 47621520 -> 37039052
$ws.Cells.Item(58, 10).Value = 2556.4443  # J58: 1758.75 -> 2556.4443
$ws.Cells.Item(58, 11).Value = 37039052  # K58: 47621520 -> 37039052
$ws.Cells.Item(58, 12).Value = 2556.4443  # L58: 1758.75 -> 2556.4443
$ws.Cells.Item(58, 13).Value = -37038849  # M58: -47621317 -> -37038849
$ws.Cells.Item(58, 14).Value = -2962.4443  # N58: -2164.75 -> -2962.4443
$ws.Cells.Item(92, 8).Value = 0  # H92: 22601 -> 0
$ws.Cells.Item(92, 10).Value = 0  # J92: 22601 -> 0
$ws.Cells.Item(92, 12).Value = 0  # L92: 22601 -> 0
$ws.Cells.Item(92, 14).ClearContents()  # N92: -27593 -> (removed)
$ws.Cells.Item(99, 8).Value = 11100  # H99: 4101.5 -> 11100
$ws.Cells.Item(99, 9).Value = 0  # I99: 1768.6666 -> 0
$ws.Cells.Item(99, 11).Value = 0  # K99: 1768.6666 -> 0
$ws.Cells.Item(99, 13).ClearContents()  # M99: -270.6666 -> (removed)
$ws.Cells.Item(126, 8).Value = 11100  # H126: 4101.5 -> 11100
$ws.Cells.Item(126, 9).Value = 0  # I126: 1768.6666 -> 0
$ws.Cells.Item(126, 11).Value = 0  # K126: 5305.9998 -> 0
$ws.Cells.Item(126, 13).ClearContents()  # M126: -2835.9998 -> (removed)
$ws.Cells.Item(132, 8).Value = 52759.5  # H132: 29838.666 -> 52759.5
$ws.Cells.Item(132, 9).Value = 2599.4285  # I132: 1928.45 -> 2599.4285
$ws.Cells.Item(132, 10).Value = 169799.67  # J132: 64726.438 -> 169799.67
$ws.Cells.Item(132, 11).Value = 7798.2855  # K132: 5785.35 -> 7798.2855
$ws.Cells.Item(132, 12).Value = 509399.01  # L132: 194179.314 -> 509399.01
$ws.Cells.Item(132, 13).Value = -5268.2855  # M132: -3255.35 -> -5268.2855
$ws.Cells.Item(132, 14).Value = -514459.01  # N132: -199239.314 -> -514459.01
$ws.Cells.Item(134, 8).Value = 29254.82  # H134: 27112.762 -> 29254.82
$ws.Cells.Item(134, 9).Value = 1779.0454  # I134: 1475.48 -> 1779.0454
$ws.Cells.Item(134, 10).Value = 64811.707  # J134: 64814.65 -> 64811.707
$ws.Cells.Item(134, 11).Value = 5337.1362  # K134: 4426.440000000001 -> 5337.1362
$ws.Cells.Item(134, 12).Value = 194435.121  # L134: 194443.95 -> 194435.121
$ws.Cells.Item(134, 13).Value = -2802.1362  # M134: -1891.440000000001 -> -2802.1362
$ws.Cells.Item(134, 14).Value = -199505.121  # N134: -199513.95 -> -199505.121
$ws.Cells.Item(136, 8).Value = 27779928  # H136: 27029192 -> 27779928
$ws.Cells.Item(136, 9).Value = 37039052  # I136: 47621520 -> 37039052
$ws.Cells.Item(136, 10).Value = 2556.4443  # J136: 1758.75 -> 2556.4443
$ws.Cells.Item(136, 11).Value = 111117156  # K136: 142864560 -> 111117156
$ws.Cells.Item(136, 12).Value = 7669.3329  # L136: 5276.25 -> 7669.3329
$ws.Cells.Item(136, 13).Value = -111114606  # M136: -142862010 -> -111114606
$ws.Cells.Item(136, 14).Value = -12769.3329  # N136: -10376.25 -> -12769.3329

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5357.615  # H3: 5415.5625 -> 5357.615
$ws.Cells.Item(3, 9).Value = 5387.4165  # I3: 5554.0835 -> 5387.4165
$ws.Cells.Item(3, 11).Value = 16162.2495  # K3: 16662.2505 -> 16162.2495
$ws.Cells.Item(3, 13).Value = -16050.2495  # M3: -16550.2505 -> -16050.2495
$ws.Cells.Item(128, 8).Value = 869250  # H128: 498500 -> 869250
$ws.Cells.Item(128, 9).Value = 869250  # I128: 498500 -> 869250
$ws.Cells.Item(128, 11).Value = 2607750  # K128: 1495500 -> 2607750
$ws.Cells.Item(128, 13).Value = -2602770  # M128: -1490520 -> -2602770
$ws.Cells.Item(131, 8).Value = 1008.3125  # H131: 1025.1791 -> 1008.3125
$ws.Cells.Item(131, 10).Value = 1107.3091  # J131: 1121.6724 -> 1107.3091
$ws.Cells.Item(131, 12).Value = 3321.9273  # L131: 3365.0172 -> 3321.9273
$ws.Cells.Item(131, 14).Value = -13401.9273  # N131: -13445.0172 -> -13401.9273
$ws.Cells.Item(137, 8).Value = 51577.5  # H137: 31277.65 -> 51577.5
$ws.Cells.Item(137, 9).Value = 982.5  # I137: 788.5714 -> 982.5
$ws.Cells.Item(137, 10).Value = 76875  # J137: 47694.848 -> 76875
$ws.Cells.Item(137, 11).Value = 2947.5  # K137: 2365.7142 -> 2947.5
$ws.Cells.Item(137, 12).Value = 230625  # L137: 143084.544 -> 230625
$ws.Cells.Item(137, 13).Value = 2152.5  # M137: 2734.2858 -> 2152.5
$ws.Cells.Item(137, 14).Value = -240825  # N137: -153284.544 -> -240825

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 200  # H102: 1187.2632 -> 200
$ws.Cells.Item(102, 9).Value = 200  # I102: 1072.6666 -> 200
$ws.Cells.Item(102, 10).Value = 0  # J102: 1617 -> 0
$ws.Cells.Item(102, 11).Value = 200  # K102: 1072.6666 -> 200
$ws.Cells.Item(102, 12).Value = 0  # L102: 1617 -> 0
$ws.Cells.Item(102, 13).Value = 1422  # M102: 549.3334 -> 1422
$ws.Cells.Item(102, 14).ClearContents()  # N102: -4861 -> (removed)
$ws.Cells.Item(111, 8).Value = 29800  # H111: 0 -> 29800
$ws.Cells.Item(111, 10).Value = 29800  # J111: 0 -> 29800
$ws.Cells.Item(111, 12).Value = 29800  # L111: 0 -> 29800
$ws.Cells.Item(111, 14).Value = -35934  # N111: None -> -35934

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5380.067  # H7: 6107.077 -> 5380.067
$ws.Cells.Item(7, 9).Value = 3182  # I7: 3581.0908 -> 3182
$ws.Cells.Item(7, 10).Value = 11424.75  # J7: 20000 -> 11424.75
$ws.Cells.Item(7, 11).Value = 3182  # K7: 3581.0908 -> 3182
$ws.Cells.Item(7, 12).Value = 11424.75  # L7: 20000 -> 11424.75
$ws.Cells.Item(7, 13).Value = -3070  # M7: -3469.0908 -> -3070
$ws.Cells.Item(7, 14).Value = -11648.75  # N7: -20224 -> -11648.75
$ws.Cells.Item(22, 8).Value = 668.0769  # H22: 694.1539 -> 668.0769
$ws.Cells.Item(22, 9).Value = 436.66666  # I22: 443.2143 -> 436.66666
$ws.Cells.Item(22, 10).Value = 770.9259  # J22: 834.6799999999999 -> 770.9259
$ws.Cells.Item(22, 11).Value = 436.66666  # K22: 443.2143 -> 436.66666
$ws.Cells.Item(22, 12).Value = 770.9259  # L22: 834.6799999999999 -> 770.9259
$ws.Cells.Item(22, 13).Value = -141.66666  # M22: -148.2143 -> -141.66666
$ws.Cells.Item(22, 14).Value = -1360.9259  # N22: -1424.68 -> -1360.9259
$ws.Cells.Item(27, 8).Value = 668.0769  # H27: 694.1539 -> 668.0769
$ws.Cells.Item(27, 9).Value = 436.66666  # I27: 443.2143 -> 436.66666
$ws.Cells.Item(27, 10).Value = 770.9259  # J27: 834.6799999999999 -> 770.9259
$ws.Cells.Item(27, 11).Value = 436.66666  # K27: 443.2143 -> 436.66666
$ws.Cells.Item(27, 12).Value = 770.9259  # L27: 834.6799999999999 -> 770.9259
$ws.Cells.Item(27, 13).Value = -329.66666  # M27: -336.2143 -> -329.66666
$ws.Cells.Item(27, 14).Value = -984.9259  # N27: -1048.68 -> -984.9259
$ws.Cells.Item(43, 8).Value = 10000  # H43: 0 -> 10000
$ws.Cells.Item(43, 10).Value = 10000  # J43: 0 -> 10000
$ws.Cells.Item(43, 12).Value = 10000  # L43: 0 -> 10000
$ws.Cells.Item(43, 14).Value = -10386  # N43: None -> -10386
$ws.Cells.Item(100, 8).Value = 1693.5186  # H100: 1795.15 -> 1693.5186
$ws.Cells.Item(100, 9).Value = 1620.4546  # I100: 1960.6 -> 1620.4546
$ws.Cells.Item(100, 10).Value = 1743.75  # J100: 1740 -> 1743.75
$ws.Cells.Item(100, 11).Value = 1620.4546  # K100: 1960.6 -> 1620.4546
$ws.Cells.Item(100, 12).Value = 1743.75  # L100: 1740 -> 1743.75
$ws.Cells.Item(100, 13).Value = -1079.4546  # M100: -1419.6 -> -1079.4546
$ws.Cells.Item(100, 14).Value = -2825.75  # N100: -2822 -> -2825.75
$ws.Cells.Item(126, 8).Value = 5380.067  # H126: 6107.077 -> 5380.067
$ws.Cells.Item(126, 9).Value = 3182  # I126: 3581.0908 -> 3182
$ws.Cells.Item(126, 10).Value = 11424.75  # J126: 20000 -> 11424.75
$ws.Cells.Item(126, 11).Value = 9546  # K126: 10743.2724 -> 9546
$ws.Cells.Item(126, 12).Value = 34274.25  # L126: 60000 -> 34274.25
$ws.Cells.Item(126, 13).Value = -7076  # M126: -8273.2724 -> -7076
$ws.Cells.Item(126, 14).Value = -39214.25  # N126: -64940 -> -39214.25
$ws.Cells.Item(134, 8).Value = 0  # H134: 59699.5 -> 0
$ws.Cells.Item(134, 10).Value = 0  # J134: 59699.5 -> 0
$ws.Cells.Item(134, 12).Value = 0  # L134: 59699.5 -> 0
$ws.Cells.Item(134, 14).ClearContents()  # N134: -69839.5 -> (removed)
$ws.Cells.Item(137, 8).Value = 30072.5  # H137: 31650 -> 30072.5
$ws.Cells.Item(137, 9).Value = 20390  # I137: 0 -> 20390
$ws.Cells.Item(137, 10).Value = 33300  # J137: 31650 -> 33300
$ws.Cells.Item(137, 11).Value = 20390  # K137: 0 -> 20390
$ws.Cells.Item(137, 12).Value = 33300  # L137: 31650 -> 33300
$ws.Cells.Item(137, 13).Value = -15290  # M137: None -> -15290
$ws.Cells.Item(137, 14).Value = -43500  # N137: -41850 -> -43500
$ws.Cells.Item(139, 8).Value = 49041.668  # H139: 47800 -> 49041.668
$ws.Cells.Item(139, 9).Value = 0  # I139: 39000 -> 0
$ws.Cells.Item(139, 10).Value = 49041.668  # J139: 50000 -> 49041.668
$ws.Cells.Item(139, 11).Value = 0  # K139: 39000 -> 0
$ws.Cells.Item(139, 12).Value = 49041.668  # L139: 50000 -> 49041.668
$ws.Cells.Item(139, 13).ClearContents()  # M139: -33860 -> (removed)
$ws.Cells.Item(139, 14).Value = -59321.668  # N139: -60280 -> -59321.668
$ws.Cells.Item(141, 8).Value = 65155.715  # H141: 59540.91 -> 65155.715
$ws.Cells.Item(141, 10).Value = 65155.715  # J141: 59540.91 -> 65155.715
$ws.Cells.Item(141, 12).Value = 65155.715  # L141: 59540.91 -> 65155.715
$ws.Cells.Item(141, 14).Value = -75515.715  # N141: -69900.91 -> -75515.715

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 55133.1  # H132: 41606.902 -> 55133.1
$ws.Cells.Item(132, 9).Value = 35352.75  # I132: 35355.188 -> 35352.75
$ws.Cells.Item(132, 10).Value = 145557.58  # J132: 51609.65 -> 145557.58
$ws.Cells.Item(132, 11).Value = 106058.25  # K132: 106065.564 -> 106058.25
$ws.Cells.Item(132, 12).Value = 436672.74  # L132: 154828.95 -> 436672.74
$ws.Cells.Item(132, 13).Value = -103528.25  # M132: -103535.564 -> -103528.25
$ws.Cells.Item(132, 14).Value = -441732.74  # N132: -159888.95 -> -441732.74
$ws.Cells.Item(136, 8).Value = 25867.549  # H136: 29457.709 -> 25867.549
$ws.Cells.Item(136, 9).Value = 18568.84  # I136: 22134.148 -> 18568.84
$ws.Cells.Item(136, 10).Value = 41587.848  # J136: 43226 -> 41587.848
$ws.Cells.Item(136, 11).Value = 55706.52  # K136: 66402.444 -> 55706.52
$ws.Cells.Item(136, 12).Value = 124763.544  # L136: 129678 -> 124763.544
$ws.Cells.Item(136, 13).Value = -53156.52  # M136: -63852.444 -> -53156.52
$ws.Cells.Item(136, 14).Value = -134863.544  # N136: -134778 -> -134863.544
